$d = $word.ActiveDocument

# 1) Title paragraph: merge "Cover Latter of " + "Front-End Developer" into a single run.
$d.Content.Find.Execute(
    "Cover Latter of Front-End Developer", $false, $false, $false, $false, $false,
    $true, 1, $false, "Cover Latter of Front-End Developer", 2)

# 2) Phone numbers: merge "9916183187" + ", 9113852294" into a single run.
$d.Content.Find.Execute(
    "9916183187, 9113852294", $false, $false, $false, $false, $false,
    $true, 1, $false, "9916183187, 9113852294", 2)

# 3) Opening paragraph: merge all the split runs ("...User ", "Experience", ", along with my
#    ability to Angular, ", "react", " make me...") into one contiguous run.
$d.Content.Find.Execute(
    "I am writing in application for the front-end developer position at [company name] as advertised in Naukri Portal. My Experience in the field of Frontend programming and User Experience, along with my ability to Angular, react make me the perfect candidate for the job. I know that I would be a valuable addition to the team at [company name].",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "I am writing in application for the front-end developer position at [company name] as advertised in Naukri Portal. My Experience in the field of Frontend programming and User Experience, along with my ability to Angular, react make me the perfect candidate for the job. I know that I would be a valuable addition to the team at [company name].",
    2)

# 4) Split the "if applicable]. Feel free to contact me directly should you require anything
#    further." paragraph: keep "if applicable]. " here, move the rest to the start of the
#    "Thank you for your time..." paragraph (after the page break), separated by a space.
$d.Content.Find.Execute(
    "if applicable]. Feel free to contact me directly should you require anything further.",
    $false, $false, $false, $false, $false,
    $true, 1, $false, "if applicable]. ", 2)

$d.Content.Find.Execute(
    "Thank you for your time and consideration. I look forward to hearing from you.",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "Feel free to contact me directly should you require anything further. Thank you for your time and consideration. I look forward to hearing from you.",
    2)
